$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 155, shifting existing rows 155..245 down to 156..246
$ws.Rows.Item(155).Insert()

# Populate the newly inserted row 155 with the new record's data
$ws.Range("A155").Value = 10
$ws.Range("B155").Value = 'Vega Modelo de Temuco'
$ws.Range("C155").Value = 'La Araucanía'
$ws.Range("D155").Value = 45236
$ws.Range("E155").Value = 9
$ws.Range("F155").Value = 'Fruta'
$ws.Range("G155").Value = 100107
$ws.Range("H155").Value = 'Otros'
$ws.Range("I155").Value = 100107002
$ws.Range("J155").Value = 'Chirimoya'
$ws.Range("K155").Value = 'Cultivar IV Región'
$ws.Range("L155").Value = 'Primera'
$ws.Range("M155").Value = 65
$ws.Range("N155").Value = 2500
$ws.Range("O155").Value = 2500
$ws.Range("P155").Value = 2500
$ws.Range("Q155").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R155").Value = 'Provincia del Elquí'
$ws.Range("S155").Value = 2500
$ws.Range("T155").Value = 1
